# Household()'s Residence() -> add Tenure Pref / Tenure columns to "households" sheet,
# fix a data value, tweak column width, and update sheet selections/active tab.

$wb = $excel.ActiveWorkbook

$households = $wb.Worksheets.Item("households")
$housingStock = $wb.Worksheets.Item("housing_stock")
$financialCapital = $wb.Worksheets.Item("financial_capital")

# --- households sheet: new "Tenure Pref" (N) and "Tenure" (O) columns ---
# Column header + values are written in this specific order so that the
# workbook's shared-string table grows in the same order as the target:
# "Tenure Pref", "Rent", "Own", "Tenure".
$households.Range("N1").Value = "Tenure Pref"
$households.Range("O2").Value = "Rent"
$households.Range("O3").Value = "Own"
$households.Range("O4").Value = "Rent"
$households.Range("O5").Value = "Own"
$households.Range("O1").Value = "Tenure"
$households.Range("N2").Value = "Own"
$households.Range("N3").Value = "Own"
$households.Range("N4").Value = "Rent"
$households.Range("N5").Value = "Own"

# Data correction: Bruce's Bedrooms count 6 -> 5
$households.Range("H3").Value = 5

# New column N gets a custom width
$households.Columns.Item(13).ColumnWidth = 16

# --- housing_stock sheet: move the remembered selection ---
$housingStock.Range("K30").Select() | Out-Null

# --- households sheet becomes the active/visible tab ---
$households.Activate() | Out-Null

$wb.Save()
